$d = $word.ActiveDocument

# 1. Remove the trailing space at the end of the "second page" description paragraph.
$d.Content.Find.Execute("features that affect eligibility. ", $true, $false, $false, $false, $false, $true, 1, $false, "features that affect eligibility.", 2) | Out-Null

# 2. Insert a new paragraph (with text) and two blank paragraphs right after it,
#    pushing the "Figure 3" screenshot further down the page.
$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs(10)
$p10.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs(11)
$p11.Range.InsertParagraphAfter()

$d.Paragraphs(10).Range.Text = "The newly created records along with their predictions will be saved to a csv file for further analysis by the business."

$emptyParaXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(11).Range.InsertXML($emptyParaXml)
$d.Paragraphs(12).Range.InsertXML($emptyParaXml)

# 3. The "Figure 3" screenshot paragraph now begins a new page, so give it the
#    rendered-page-break marker it picked up.
$metricsShotXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="1C3200CD" w14:textId="77777777" w:rsidR="009E0668" w:rsidRDefault="009E0668" w:rsidP="009E0668"><w:r w:rsidRPr="00AC583C"><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="14CACEEF" wp14:editId="0DCB53B2"><wp:extent cx="5189220" cy="2805284"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1719597428" name="Picture 1" descr="A screenshot of a computer&#xA;&#xA;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1719597428" name="Picture 1" descr="A screenshot of a computer&#xA;&#xA;Description automatically generated"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5194774" cy="2808287"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(13).Range.InsertXML($metricsShotXml)

# 4. The "Figure 4" screenshot paragraph no longer starts a page, so drop its
#    rendered-page-break marker.
$devTeamShotXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="21AD7CA0" w14:textId="77777777" w:rsidR="009E0668" w:rsidRDefault="009E0668" w:rsidP="009E0668"><w:r w:rsidRPr="005C5D48"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="70DE13E8" wp14:editId="12AE9C93"><wp:extent cx="5943600" cy="2506345"/><wp:effectExtent l="0" t="0" r="0" b="8255"/><wp:docPr id="38632522" name="Picture 1" descr="A screenshot of a computer&#xA;&#xA;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="38632522" name="Picture 1" descr="A screenshot of a computer&#xA;&#xA;Description automatically generated"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId7"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="2506345"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(16).Range.InsertXML($devTeamShotXml)

Write-Output "done"
